$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1417620482754679
$ws.Range("C2").Value = 0.5996277706543274
$ws.Range("D2").Value = 0.8169135268031753
$ws.Range("E2").Value = 0.9038326873947276
$ws.Range("F2").Value = 0.9017087228643812
$ws.Range("G2").Value = 50

$ws.Range("B3").Value = 0.09584318641535604
$ws.Range("C3").Value = 0.657502611948801
$ws.Range("D3").Value = 1.03928723612716
$ws.Range("E3").Value = 1.019454381582207
$ws.Range("F3").Value = 1.025456856189602
$ws.Range("G3").Value = 49

$ws.Range("B4").Value = 0.126342760679196
$ws.Range("C4").Value = 0.5716138419826224
$ws.Range("D4").Value = 0.7986488544558821
$ws.Range("E4").Value = 0.8936715584910835
$ws.Range("F4").Value = 0.894057751257528
$ws.Range("G4").Value = 48

$ws.Range("B5").Value = 0.132571133247816
$ws.Range("C5").Value = 0.6727079862838818
$ws.Range("D5").Value = 1.180166600671597
$ws.Range("E5").Value = 1.086354730588309
$ws.Range("F5").Value = 1.089892299018441
$ws.Range("G5").Value = 47

$ws.Range("B6").Value = 0.1115124553593621
$ws.Range("C6").Value = 0.6125314060269084
$ws.Range("D6").Value = 0.941813663750834
$ws.Range("E6").Value = 0.9704708464198365
$ws.Range("F6").Value = 0.9746955907510105
$ws.Range("G6").Value = 46

$ws.Range("B7").Value = 0.0971729415500782
$ws.Range("C7").Value = 0.5646146643790109
$ws.Range("D7").Value = 0.6953729178386839
$ws.Range("E7").Value = 0.8338902312886775
$ws.Range("F7").Value = 0.8406640857472294
$ws.Range("G7").Value = 34

$ws.Range("B8").Value = 0.1456103435239261
$ws.Range("C8").Value = 0.5758738547837706
$ws.Range("D8").Value = 0.7130976704205472
$ws.Range("E8").Value = 0.8444511059975865
$ws.Range("F8").Value = 0.8446993703979081
$ws.Range("G8").Value = 33

$ws.Range("B9").Value = -0.01600337040817491
$ws.Range("C9").Value = 0.5629050319730056
$ws.Range("D9").Value = 0.5311176513683245
$ws.Range("E9").Value = 0.7287781907880644
$ws.Range("F9").Value = 0.7524973840956725
$ws.Range("G9").Value = 16

$ws.Range("B10").Value = -0.08801024940629837
$ws.Range("C10").Value = 0.472497539938865
$ws.Range("D10").Value = 0.4690803580242086
$ws.Range("E10").Value = 0.6848944137779257
$ws.Range("F10").Value = 0.7159566669255664
$ws.Range("G10").Value = 10

$ws.Range("B11").Value = 0.1495946378431339
$ws.Range("C11").Value = 0.4169860915672885
$ws.Range("D11").Value = 0.2389570720534238
$ws.Range("E11").Value = 0.4888323557759078
$ws.Range("F11").Value = 0.5203106240290569
